$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (P1-001): append a new history line to the Notas/Historial note ---
$e3 = $ws.Cells.Item(3, 5)
$e3.Value = $e3.Value() + "`n2/6: General check, no problems other than copper areas in mounting holes."
$ws.Rows.Item(3).RowHeight = 68

# --- Row 9 (new): P1-007 ---
$ws.Cells.Item(9, 1).Value = "P1-007"
$ws.Cells.Item(9, 2).Value = 38749
$ws.Cells.Item(9, 2).NumberFormat = "mmm-yy"
$ws.Cells.Item(9, 3).Value = "x"
$ws.Cells.Item(9, 5).Value = "2/6: Micro falso? No se pudo programar a la primera."

# --- Row 10 (new): P1-008 ---
$ws.Cells.Item(10, 1).Value = "P1-008"
$ws.Cells.Item(10, 2).Value = 38749
$ws.Cells.Item(10, 2).NumberFormat = "mmm-yy"
$ws.Cells.Item(10, 3).Value = "pre-1.0.6"
$ws.Cells.Item(10, 5).Value = "2/6: Ensamble, probada."

# --- Selection ends at E14 ---
$ws.Range("E14").Select() | Out-Null
